$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank row at position 17. This pushes the old row 17 (empty) and
#    everything below it down by one: old 18 (header) -> 19, old 19 -> 20, old 20 -> 21,
#    old 21 -> 22, old 22 -> 23.
$ws.Rows.Item(17).Insert()

# 2. New row 17: mirrors the C:H structure of row 16 with a new base_channel=256 case,
#    tagged with the git version marker in column N.
$ws.Range("C17").Value = 256
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 128
$ws.Range("F17").Value = 112
$ws.Range("G17").Value = 112
$ws.Range("H17").Value = 3
$ws.Range("N17").Value = "bb29928"

# 3. Row 16 gains a note in column A, the measured timings in J:M, and the git version tag in N.
$ws.Range("A16").Value = "改回来了，加入v3"
$ws.Range("J16").Value = 0.0045
$ws.Range("K16").Value = 0.0213
$ws.Range("L16").Value = 0.0134
$ws.Range("M16").Value = 0.0126
$ws.Range("N16").Value = "bb29928"

# 4. Row 13 header gains the new FLKv3 column.
$ws.Range("M13").Value = "FLKv3"

# 5. Row 5 banner gains the git version label, mirroring the "one conv time(s)" banner.
$ws.Range("N5").Value = "git version"

# 6. Row 19 (the base_channel header, shifted down from the old row 18) no longer carries
#    the git-version header cell that used to live in column N.
$ws.Range("N19").ClearContents()

# The old annotation that used to live in column A (shifted down to row 20 by the insert)
# is removed; it is superseded by the new note added to row 16.
$ws.Range("A20").ClearContents()

# 7. The four data rows that shifted down (20:23) each gain the git version tag in column N.
$ws.Range("N20").Value = "6a9aec0"
$ws.Range("N21").Value = "6a9aec0"
$ws.Range("N22").Value = "6a9aec0"
$ws.Range("N23").Value = "6a9aec0"

# 8. Update the sheet dimension + selection bookkeeping to match the new extent.
$ws.Range("N23").Select()
